$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.072.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3858"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07849"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9609"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.891.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.674"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.887"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06842"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009925"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.102.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.290"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.106.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.754"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.968"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9416"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09244"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.265"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.320"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.331"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02106"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.721"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5593"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.899"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07327"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("E46").Value = "  -8.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.110"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.833"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.020"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
